$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated (bugfixed) values for existing rows 2-9 ---
$data = @{
    2  = @{ B = 0.1805204829727294;   C = 0.9134701371387803; D = 2.313862988505206;  E = 1.521138714419302;  F = 1.567404881478029 }
    3  = @{ B = 0.1654443250838863;   C = 0.8321534822147496; D = 2.430194309093069;   E = 1.55890805023679;   F = 1.613399420919432 }
    4  = @{ B = 0.3342101850538108;   C = 0.8990974440469651; D = 3.047122944578856;   E = 1.745601026746621;  F = 1.789492512149283 }
    5  = @{ B = 0.349057914304527;    C = 0.9234061813053621; D = 2.965460108284523;   E = 1.722051134050474;  F = 1.768609778560688 }
    6  = @{ B = 0.4342163702162936;   C = 1.060274838210373;  D = 3.344243897922033;   E = 1.828727398471963;  F = 1.872520595275618 }
    7  = @{ B = -0.06951841543571383; C = 0.5087111458309537; D = 0.3139214969299044;  E = 0.560286977298156;  F = 0.5896819250247234; G = 9 }
    8  = @{ B = -0.04657217614193387; C = 0.7130288563221177; D = 0.6707611900940872;  E = 0.8190001160525481; F = 0.8957179617514601; G = 6 }
    9  = @{ B = -0.3788153344042176;  C = 0.3788153344042176; D = 0.2155603798649167;  E = 0.4642848046887995; F = 0.3287688906020555; G = 3 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# --- Add new row 10 for a new simulated quarter ("Q8") ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats - copy the same border/bold/center style as A2:A9
$excel.CutCopyMode = $false

$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.410128014204804
$ws.Range("C10").Value = 0.410128014204804
$ws.Range("D10").Value = 0.168204988035576
$ws.Range("E10").Value = 0.410128014204804
$ws.Range("G10").Value = 1
